$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: Elapsed Time / CPU
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Copy the existing header formatting (bold, border, centered) onto the new headers
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

# Updated metrics for row 2
$ws.Range("B2").Value = 0.2139593575861684
$ws.Range("C2").Value = 0.9841919032133989
$ws.Range("D2").Value = 0.3608222963467753
$ws.Range("F2").Value = "Pipeline(steps=[('model', AdaBoostRegressor(n_estimators=150))])"

# New elapsed time / cpu values for row 2
$ws.Range("G2").Value = 0.1434780816666413
$ws.Range("H2").Value = 0.992
